# "March 24 update 3" — add three new trailing columns (renewd / PlanID /
# iteration) to the bldg sheet and fill them in for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (M1:O1), mirroring the existing header row (B1:L1).
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Give the new headers the same look (bold font, border, centered) as the
# rest of row 1 by copying the formatting from the last existing header.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new columns on every data row (2-45).
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"   # M: renewd
    $ws.Cells.Item($r, 14).Value = 20181295  # N: PlanID
    $ws.Cells.Item($r, 15).Value = 7         # O: iteration
}
